$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct the D6 (AMD buy) units to a round number ---
$ws.Range("D6").Value = 21

# --- Add the "cost" column (F) with a units*price formula, filled down ---
$ws.Range("F1").Value = "cost"
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F2").Formula = "=D2*E2"
$ws.Range("F3:F12").Formula = "=D3*E3"

# --- Add a new transaction row: selling the AMD position ---
$ws.Rows("13:13").Insert()
$ws.Range("A13").Value = 44659
$ws.Range("B13").Value = "Sell"
$ws.Range("C13").Value = "AMD"
$ws.Range("D13").Value = -21
$ws.Range("E13").Value = 101
$ws.Range("F13").Formula = "=D13*E13"

$ws.Range("D14").Select()
